$wb = $excel.ActiveWorkbook

# ---- Building sheet updates ----
$building = $wb.Worksheets.Item("Building")
$building.Range("C1").Value = "Request"
$building.Range("C3").Value = 0
$building.Range("D3").Value = 0
$building.Range("F1").Value = "Production"
$building.Range("F2").Value = "carrot"
$building.Range("F5").Value = 1

# ---- New "Natrual" sheet, inserted right after Building ----
$natrual = $wb.Worksheets.Add($null, $building)
$natrual.Name = "Natrual"

$natrual.Range("A1").Value = "Type"
$natrual.Range("B1").Value = "Name"
$natrual.Range("C1").Value = "Requests"
$natrual.Range("A1:A2").Merge()
$natrual.Range("B1:B2").Merge()
$natrual.Range("C1:D1").Merge()

$natrual.Range("C2").Value = "carrot"
$natrual.Range("D2").Value = "wood"

$natrual.Range("F1").Value = "Request"
$natrual.Range("F2").Value = "carrot"
$natrual.Range("G2").Value = "wood"

$natrual.Range("A3").Value = 0
$natrual.Range("B3").Value = "Tree"
$natrual.Range("C3").Value = 0
$natrual.Range("D3").Value = 0
$natrual.Range("F3").Value = 0
$natrual.Range("G3").Value = 0

# ---- Unit sheet updates ----
$unit = $wb.Worksheets.Item("Unit")
$unit.Range("E1").Value = "Consumes"
$unit.Range("E2").Value = "food"
$unit.Range("E3").Value = 1

# ---- Restore selections / active sheet (Building stays the active tab) ----
$unit.Range("E4").Select()
$natrual.Range("H3").Select()
$building.Range("F6").Select()

Write-Host "done"
